$p = $ppt.ActivePresentation

foreach ($idx in 15, 16) {
    $s = $p.Slides.Item($idx)
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "Subtitle 2") {
            $shp.Delete()
        }
    }
}
